# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# The "Periodo Mora" list (column E, rows 16-38) used to be ordered newest
# to oldest (2503 down to 2305). The update re-orders it oldest to newest
# (2305 up to 2503), keeping every row's other data (worker, value, etc.)
# untouched. As a side effect of that re-order, the "Valor Mora" date-style
# number in column F swaps between the first and last data row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New ascending order for the "Periodo Mora" column (rows 16 .. 38)
$periods = @(
    "2305","2306","2307","2308","2309","2310","2311","2312",
    "2401","2402","2403","2404","2405","2406","2407","2408","2409","2410","2411","2412",
    "2501","2502","2503"
)

$firstRow = 16
for ($i = 0; $i -lt $periods.Length; $i++) {
    $row = $firstRow + $i
    $ws.Cells.Item($row, 5).Value = $periods[$i]
}

# Column F ("Valor Mora") values swap between the first row (16) and the
# last row (38) of the table.
$ws.Range("F16").Value = 71240
$ws.Range("F38").Value = 42744
